$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet from "Sheet1" to "raw data"
$ws.Name = "raw data"

# Fix header text in D1: "Answer_relevance" -> "Answer relevance"
$ws.Range("D1").Value = "Answer relevance"

# Update the scroll position / active cell selection for the sheet view
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("D2").Select() | Out-Null
